$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59. Excel shifts rows 59..160 down to 60..161,
# and the previously-last row (160) gets duplicated into the new row 161
# (matching the target dimension A1:T161).
$ws.Rows.Item(59).Insert()

# Fill the newly-inserted (now-empty) row 59 with the new data record.
$ws.Cells.Item(59, 1).Value = 8
$ws.Cells.Item(59, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44581
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100103
$ws.Cells.Item(59, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(59, 9).Value = 100103002
$ws.Cells.Item(59, 10).Value = "Ciruela"
$ws.Cells.Item(59, 11).Value = "Black Amber"
$ws.Cells.Item(59, 12).Value = "Segunda"
$ws.Cells.Item(59, 13).Value = 20
$ws.Cells.Item(59, 14).Value = 225000
$ws.Cells.Item(59, 15).Value = 230000
$ws.Cells.Item(59, 16).Value = 227500
$ws.Cells.Item(59, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(59, 18).Value = "Región Metropolitana"
$ws.Cells.Item(59, 19).Value = 506
$ws.Cells.Item(59, 20).Value = 450
